# fix(publipostage): Correct status, status name, status label ...
#
# - Drop the "statut_label" column (B): its (vert/noir) values duplicated
#   information already carried by "statut"/"statut_name".
# - Drop the "results_1y" / "results_3y" / "results" boolean columns: no
#   longer part of the mail-merge output.
# - Re-encode "statut" as the numeric code (1 / 4) instead of the emoji
#   (✅ / ⚠️), and prefix "statut_name" with that same numeric code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove obsolete columns -------------------------------------------
# "statut_label" (vert/noir) is no longer needed.
$ws.Columns("B").Delete()

# After the delete above, "results_1y"/"results_3y"/"results" sit at I:K
# (they used to be J:L, shifted left by the column-B removal).
$ws.Range("I1:K1").EntireColumn.Delete()

# --- Recode "statut" (A) from emoji to numeric status code --------------
# Route the new value through a TEXT() formula + paste-values so the cell
# keeps a genuine text type ("1"/"4") instead of Excel auto-coercing a
# bare numeric-looking string into a number.
$ws.Range("A2").Formula = '=TEXT(1,"0")'
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)

$ws.Range("A3").Formula = '=TEXT(4,"0")'
$ws.Range("A3").Copy()
$ws.Range("A3").PasteSpecial(-4163)

$ws.Range("A4").Formula = '=TEXT(4,"0")'
$ws.Range("A4").Copy()
$ws.Range("A4").PasteSpecial(-4163)

# --- Update "statut_name" (B) text to include the numeric code prefix ---
$ws.Range("B2").Value = "1: résultats postés ou publiés dans les 12 mois"
$ws.Range("B3").Value = "4: pas de résultats postés ni publiés"
$ws.Range("B4").Value = "4: pas de résultats postés ni publiés"
